$d = $word.ActiveDocument

# Header(2) holds the "BTec_Logo-Orange" inline picture (header1.xml):
#   name: image1.jpg -> image2.jpg
$d.Sections(1).Headers(2).Range.InlineShapes(1).Name = "image2.jpg"

# Footer(1) holds the Pearson logo inline picture (footer1.xml, docPr id=3):
#   name: image2.png -> image1.png
$d.Sections(1).Footers(1).Range.InlineShapes(1).Name = "image1.png"

# Footer(2) holds the Pearson logo inline picture (footer2.xml, docPr id=2):
#   name: image2.png -> image1.png
$d.Sections(1).Footers(2).Range.InlineShapes(1).Name = "image1.png"
